$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 18757550
$ws.Range("I74").Value = 37506600
$ws.Range("K74").Value = 37506600
$ws.Range("M74").Value = -37505664

$ws.Range("H77").Value = 18757550
$ws.Range("I77").Value = 37506600
$ws.Range("K77").Value = 187533000
$ws.Range("M77").Value = -187528320

$ws.Range("H92").Value = 938.2632
$ws.Range("I92").Value = 823.44446
$ws.Range("K92").Value = 823.44446
$ws.Range("M92").Value = 424.55554

$ws.Range("H100").Value = 4372
$ws.Range("I100").Value = 4267.5
$ws.Range("K100").Value = 4267.5
$ws.Range("M100").Value = -3726.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 972.9583
$ws.Range("I2").Value = 1020.619
$ws.Range("K2").Value = 1020.619
$ws.Range("M2").Value = -907.619

$ws.Range("H32").Value = 6812.507
$ws.Range("I32").Value = 3552.0938
$ws.Range("K32").Value = 3552.0938
$ws.Range("M32").Value = -3265.0938

$ws.Range("H74").Value = 21745930
$ws.Range("I74").Value = 38463844
$ws.Range("J74").Value = 12641
$ws.Range("K74").Value = 38463844
$ws.Range("L74").Value = 12641
$ws.Range("M74").Value = -38462970
$ws.Range("N74").Value = -14389

$ws.Range("H77").Value = 21745930
$ws.Range("I77").Value = 38463844
$ws.Range("J77").Value = 12641
$ws.Range("K77").Value = 192319220
$ws.Range("L77").Value = 63205
$ws.Range("M77").Value = -192314852
$ws.Range("N77").Value = -71941

$ws.Range("H116").Value = 972.9583
$ws.Range("I116").Value = 1020.619
$ws.Range("K116").Value = 1020.619
$ws.Range("M116").Value = 1273.381

$ws.Range("H132").Value = 4520.4443
$ws.Range("I132").Value = 2633.5715
$ws.Range("K132").Value = 7900.7145
$ws.Range("M132").Value = -5370.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 972.9583
$ws.Range("I3").Value = 1020.619
$ws.Range("K3").Value = 1020.619
$ws.Range("M3").Value = -906.619

$ws.Range("H86").Value = 4017.3076
$ws.Range("I86").Value = 4572.1113
$ws.Range("J86").Value = 2769
$ws.Range("K86").Value = 4572.1113
$ws.Range("L86").Value = 2769
$ws.Range("M86").Value = -3449.1113
$ws.Range("N86").Value = -5015

$ws.Range("H89").Value = 4017.3076
$ws.Range("I89").Value = 4572.1113
$ws.Range("J89").Value = 2769
$ws.Range("K89").Value = 22860.5565
$ws.Range("L89").Value = 13845
$ws.Range("M89").Value = -17244.5565
$ws.Range("N89").Value = -25077

$ws.Range("H99").Value = 1655.625
$ws.Range("I99").Value = 1461.4615
$ws.Range("K99").Value = 1461.4615
$ws.Range("M99").Value = 36.53850000000011

$ws.Range("H103").Value = 58739
$ws.Range("J103").Value = 58739
$ws.Range("L103").Value = 58739
$ws.Range("N103").Value = -61083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2581.8572
$ws.Range("I22").Value = 349.57144
$ws.Range("J22").Value = 4814.143
$ws.Range("K22").Value = 349.57144
$ws.Range("L22").Value = 4814.143
$ws.Range("M22").Value = 0.4285600000000045
$ws.Range("N22").Value = -5514.143

$ws.Range("H31").Value = 6459.129
$ws.Range("I31").Value = 2349.4075
$ws.Range("J31").Value = 34199.75
$ws.Range("K31").Value = 2349.4075
$ws.Range("L31").Value = 34199.75
$ws.Range("M31").Value = -2054.4075
$ws.Range("N31").Value = -34789.75

$ws.Range("H34").Value = 6459.129
$ws.Range("I34").Value = 2349.4075
$ws.Range("J34").Value = 34199.75
$ws.Range("K34").Value = 2349.4075
$ws.Range("L34").Value = 34199.75
$ws.Range("M34").Value = -2147.4075
$ws.Range("N34").Value = -34603.75

$ws.Range("H50").Value = 64568.855
$ws.Range("I50").Value = 40663.332
$ws.Range("J50").Value = 82498
$ws.Range("K50").Value = 40663.332
$ws.Range("L50").Value = 82498
$ws.Range("M50").Value = -40038.332
$ws.Range("N50").Value = -83748

$ws.Range("H58").Value = 3796.6553
$ws.Range("I58").Value = 2753
$ws.Range("J58").Value = 5779.6
$ws.Range("K58").Value = 2753
$ws.Range("L58").Value = 5779.6
$ws.Range("M58").Value = -2550
$ws.Range("N58").Value = -6185.6

$ws.Range("H93").Value = 9907
$ws.Range("I93").Value = 9907
$ws.Range("K93").Value = 9907
$ws.Range("M93").Value = -8035

$ws.Range("H107").Value = 1203.3572
$ws.Range("I107").Value = 986.86957
$ws.Range("J107").Value = 2199.2
$ws.Range("K107").Value = 986.86957
$ws.Range("L107").Value = 2199.2
$ws.Range("M107").Value = 933.13043
$ws.Range("N107").Value = -6039.2

$ws.Range("H122").Value = 3844.2122
$ws.Range("I122").Value = 3578.64
$ws.Range("J122").Value = 4674.125
$ws.Range("K122").Value = 10735.92
$ws.Range("L122").Value = 14022.375
$ws.Range("M122").Value = -8285.92
$ws.Range("N122").Value = -18922.375

$ws.Range("H132").Value = 4667.222
$ws.Range("I132").Value = 2866.1428
$ws.Range("J132").Value = 10971
$ws.Range("K132").Value = 8598.428400000001
$ws.Range("L132").Value = 32913
$ws.Range("M132").Value = -6068.428400000001
$ws.Range("N132").Value = -37973

$ws.Range("H136").Value = 3796.6553
$ws.Range("I136").Value = 2753
$ws.Range("J136").Value = 5779.6
$ws.Range("K136").Value = 8259
$ws.Range("L136").Value = 17338.8
$ws.Range("M136").Value = -5709
$ws.Range("N136").Value = -22438.8

$ws.Range("H140").Value = 123499.2
$ws.Range("J140").Value = 123999.11
$ws.Range("L140").Value = 123999.11
$ws.Range("N140").Value = -134359.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 6064.8184
$ws.Range("J75").Value = 5971.3
$ws.Range("L75").Value = 17913.9
$ws.Range("N75").Value = -19909.9

$ws.Range("H78").Value = 6064.8184
$ws.Range("J78").Value = 5971.3
$ws.Range("L78").Value = 53741.7
$ws.Range("N78").Value = -63725.7

$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H106").Value = 17983.8
$ws.Range("J106").Value = 17983.8
$ws.Range("L106").Value = 53951.39999999999
$ws.Range("N106").Value = -55843.39999999999

$ws.Range("H132").Value = 3388.2727
$ws.Range("J132").Value = 3400.25
$ws.Range("L132").Value = 30602.25
$ws.Range("N132").Value = -35662.25

$ws.Range("H137").Value = 3560.1667
$ws.Range("I137").Value = 2622
$ws.Range("J137").Value = 4029.25
$ws.Range("K137").Value = 7866
$ws.Range("L137").Value = 12087.75
$ws.Range("M137").Value = -2766
$ws.Range("N137").Value = -22287.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1036.875
$ws.Range("I97").Value = 917.3333
$ws.Range("J97").Value = 1190.5714
$ws.Range("K97").Value = 917.3333
$ws.Range("L97").Value = 1190.5714
$ws.Range("M97").Value = -421.3333
$ws.Range("N97").Value = -2182.5714

$ws.Range("H113").Value = 2305.7
$ws.Range("J113").Value = 2799.2
$ws.Range("L113").Value = 2799.2
$ws.Range("N113").Value = -7139.2

$ws.Range("H122").Value = 4281.5625
$ws.Range("I122").Value = 1808.1538
$ws.Range("K122").Value = 5424.4614
$ws.Range("M122").Value = -2974.4614

$ws.Range("H126").Value = 3759.5
$ws.Range("I126").Value = 3513.24
$ws.Range("K126").Value = 10539.72
$ws.Range("M126").Value = -8069.719999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 67382.5
$ws.Range("J69").Value = 67382.5
$ws.Range("L69").Value = 67382.5
$ws.Range("N69").Value = -69004.5

$ws.Range("H72").Value = 67382.5
$ws.Range("J72").Value = 67382.5
$ws.Range("L72").Value = 202147.5
$ws.Range("N72").Value = -210259.5

$ws.Range("H122").Value = 5152.0454
$ws.Range("I122").Value = 4234.75
$ws.Range("J122").Value = 6252.8
$ws.Range("K122").Value = 12704.25
$ws.Range("L122").Value = 18758.4
$ws.Range("M122").Value = -10254.25
$ws.Range("N122").Value = -23658.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 26799.8
$ws.Range("J105").Value = 26799.8
$ws.Range("L105").Value = 26799.8
$ws.Range("N105").Value = -33787.8

$ws.Range("H122").Value = 2633.6667
$ws.Range("I122").Value = 2827
$ws.Range("K122").Value = 8481
$ws.Range("M122").Value = -6031

$ws.Range("H136").Value = 3455.8096
$ws.Range("I136").Value = 2661.7368
$ws.Range("K136").Value = 7985.2104
$ws.Range("M136").Value = -5435.2104
